$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new period labels for the newly appended quarters ---
$ws.Range("Q1").Value = "31/12/2023"
$ws.Range("R1").Value = "31/03/2024"
$ws.Range("S1").Value = "30/06/2024"

# Match the header formatting (bold, centered, thin border) used by the other period columns
$ws.Range("P1").Copy()
$ws.Range("Q1:S1").PasteSpecial(-4122)

# --- Data rows 2-80: new values for columns Q (31/12/2023), R (31/03/2024), S (30/06/2024) ---
$newData = @{
  2 = @(2039019.008, 2029517.056, 2064050.048)
  3 = @(1731714.048, 1769204.992, 1736721.024)
  4 = @(312182.016, 213588, 231115.008)
  5 = @(113228, 69206, 146954)
  6 = @(418484.992, 519543.008, 526096)
  7 = @(0, 0, 0)
  8 = @(0, 0, 0)
  9 = @(10377, 10986, 11637)
  10 = @(0, 0, 0)
  11 = @(877441.9840000001, 955881.9840000001, 820918.976)
  12 = @(283430.016, 229262, 287326.016)
  13 = @(0, 0, 0)
  14 = @(0, 0, 0)
  15 = @(626, 626, 644)
  16 = @(151891.008, 158567.008, 218600.992)
  17 = @(0, 0, 0)
  18 = @(0, 0, 0)
  19 = @(0, 0, 0)
  20 = @(0, 0, 0)
  21 = @(0, 0, 0)
  22 = @(1061, 1061, 813)
  23 = @(13399, 19002, 26254)
  24 = @(9415, 10987, 12936)
  25 = @(0, 0, 0)
  26 = @(2039019.008, 2029517.056, 2064050.048)
  27 = @(657486.0159999999, 622918.0159999999, 551192)
  28 = @(29071, 31469, 31495)
  29 = @(246712.992, 249534, 208263.008)
  30 = @(36202, 24198, 29038)
  31 = @(134708, 125444, 103365)
  32 = @(8250, 8217, 8217)
  33 = @(63782, 0, 0)
  34 = @(106794, 144670, 131348)
  35 = @(31966, 39386, 39466)
  36 = @(0, 0, 0)
  37 = @(735419.008, 753779.968, 757297.024)
  38 = @(252827.008, 244400.992, 274281.984)
  39 = @(2422, 2422, 2422)
  40 = @(16809, 14517, 14817)
  41 = @(7817, 16027, 16027)
  42 = @(0, 0, 0)
  43 = @(455544, 476412.992, 449748.992)
  44 = @(0, 0, 0)
  45 = @(0, 0, 0)
  46 = @(28967, 29668, 36278)
  47 = @(617146.9840000001, 623151.008, 719283.024)
  48 = @(175228, 175228, 175228)
  49 = @(-26181, -26181, -26181)
  50 = @(0, 0, 0)
  51 = @(468100, 474104, 570236.032)
  52 = @(0, 0, 0)
  53 = @(0, 0, 0)
  54 = @(0, 0, 0)
  55 = @(0, 0, 0)
  56 = @(0, 0, 0)
  57 = @($null, $null, $null)
  58 = @($null, $null, $null)
  59 = @(639620.992, 500528, 697552)
  60 = @(-429542.944, -334060, -467931.008)
  61 = @(210077.984, 166468, 229620.992)
  62 = @(117051, 0, -60801)
  63 = @(-33177, -36114, -38345)
  64 = @(0, 0, 0)
  65 = @(20829, -18373, -11463)
  66 = @(-214816, -54754, 0)
  67 = @(-132, 0, -248)
  68 = @(-289, 3823, -465)
  69 = @(8690, 10577, 7077)
  70 = @(-8979, -6754, -7542)
  71 = @($null, $null, $null)
  72 = @($null, $null, $null)
  73 = @($null, $null, $null)
  74 = @(99543.984, 61050, 118299)
  75 = @(-12888, -10669, -13393)
  76 = @(-941, -702, -1002)
  77 = @($null, $null, $null)
  78 = @($null, $null, $null)
  79 = @(-2633, -8182, -8611)
  80 = @(83081.984, 41497, 95293)
}

foreach ($r in $newData.Keys) {
  $vals = $newData[$r]
  if ($null -ne $vals[0]) { $ws.Cells.Item($r, 17).Value = $vals[0] }
  if ($null -ne $vals[1]) { $ws.Cells.Item($r, 18).Value = $vals[1] }
  if ($null -ne $vals[2]) { $ws.Cells.Item($r, 19).Value = $vals[2] }
}